$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: iaest-measure:X -> iaest-dimension:X
$ws.Range("F3").Value = "iaest-dimension:grado"
$ws.Range("I3").Value = "iaest-dimension:sexo"
$ws.Range("J3").Value = "iaest-dimension:grandes-grupos"

# Row 4: medida -> dim
$ws.Range("F4").Value = "dim"
$ws.Range("I4").Value = "dim"
$ws.Range("J4").Value = "dim"

# Row 5: xsd:string -> skos:Concept
$ws.Range("F5").Value = "skos:Concept"
$ws.Range("I5").Value = "skos:Concept"
$ws.Range("J5").Value = "skos:Concept"

# Row 6 (new row): mapping files, copy formatting from row 5 so the new
# cells share the same style as the rest of the sheet.
$ws.Range("F5").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("I5").Copy()
$ws.Range("I6").PasteSpecial(-4122)
$ws.Range("J5").Copy()
$ws.Range("J6").PasteSpecial(-4122)

$ws.Range("F6").Value = "mapping-grado.xlsx"
$ws.Range("I6").Value = "mapping-sexo.xlsx"
$ws.Range("J6").Value = "mapping-grandes-grupos.xlsx"
